$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$newRows = @(
    @("2026-02-01","16:03:18","16:00","Living Room","PRESENCE_DETECTED","Active"),
    @("2026-02-01","16:03:26","16:00","Living Room","PRESENCE_DETECTED","Active"),
    @("2026-02-01","16:03:37","16:00","Living Room","PRESENCE_DETECTED","Active"),
    @("2026-02-01","16:03:47","16:00","Living Room","PRESENCE_DETECTED","Active"),
    @("2026-02-01","16:03:58","16:00","Living Room","PRESENCE_DETECTED","Active"),
    @("2026-02-01","16:04:08","16:00","Living Room","PRESENCE_DETECTED","Active")
)

$startRow = 68
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    # Leading apostrophe forces the date-like text to stay as text instead of
    # being auto-converted into a date serial number, matching the source log's
    # plain-text "Date" column.
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
